$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values for columns A..J (numeric) and AC (numeric) are identical for the
# two new rows (6 and 7); columns K..AB stay present but blank (empty text),
# matching the existing "blank placeholder" cells used elsewhere on row 1/2
# for the other (un-run) repetition blocks.

$rows = @(6, 7)

foreach ($r in $rows) {
    $ws.Range("A$r").Value = 4 + ($r - 6)
    $ws.Range("B$r").Value = 73600
    $ws.Range("C$r").Value = 222
    $ws.Range("D$r").Value = 3
    $ws.Range("E$r").Value = 1
    $ws.Range("F$r").Value = 2
    $ws.Range("G$r").Value = 0.013513513513513514
    $ws.Range("H$r").Value = 1
    $ws.Range("I$r").Value = 0.5
    $ws.Range("J$r").Value = 0.9864864864864865

    # K..AB: blank text placeholders (no recorded run data), same shape as
    # the blank "continuation" cells already present on rows 1 and elsewhere.
    $blankCols = @("K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB")
    foreach ($col in $blankCols) {
        $ws.Range("$col$r").Formula = '=""'
    }

    $ws.Range("AC$r").Value = 0
}
